# DataOne-Registration.Sami.xlsx edit:
# "Data file and Registration test case added"
#
# On Sheet2 (the active sheet), the EmailAddress column (C) for the three
# registration test rows is updated to new test-user addresses:
#   C2: ncitester11@nih.gov -> ncitester15@nih.gov
#   C3: ncitester12@nih.gov -> ncitester16@nih.gov
#   C4: ncitester13@nih.gov -> ncitester17@nih.gov
# The FirstName/LastName columns (A, B) are unchanged.
# The selected cell also moves to C4.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the email addresses (column C) for rows 2-4.
$ws.Range("C2").Value = "ncitester15@nih.gov"
$ws.Range("C3").Value = "ncitester16@nih.gov"
$ws.Range("C4").Value = "ncitester17@nih.gov"

# The cells already carry the "Hyperlink" style with mailto: hyperlinks;
# refresh the hyperlink targets so they point at the new addresses.
$ws.Hyperlinks.Delete()
$ws.Hyperlinks.Add($ws.Range("C2"), "mailto:ncitester15@nih.gov") | Out-Null
$ws.Hyperlinks.Add($ws.Range("C3"), "mailto:ncitester16@nih.gov") | Out-Null
$ws.Hyperlinks.Add($ws.Range("C4"), "mailto:ncitester17@nih.gov") | Out-Null

# Adding the hyperlinks can introduce a slightly different (but visually
# identical) style record; re-stamp the original "Hyperlink" cell format
# (still present on the neighboring, untouched D column) to keep the
# formatting exactly as it was.
$ws.Range("D2").Copy() | Out-Null
$ws.Range("C2").PasteSpecial(-4122) | Out-Null
$ws.Range("D3").Copy() | Out-Null
$ws.Range("C3").PasteSpecial(-4122) | Out-Null
$ws.Range("D2").Copy() | Out-Null
$ws.Range("C4").PasteSpecial(-4122) | Out-Null

# Match the saved selection state (active cell C4).
$ws.Range("C4").Select() | Out-Null
